$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove entire job entries / achievement sub-sections that were
#    dropped in the edit. Delete from the bottom of the document
#    upward so earlier paragraph indices stay valid.
# ------------------------------------------------------------------

# "Big Data & Performance Engineering" + "Technical Leadership & Integration"
# sub-sections under KEY ACHIEVEMENTS AND IMPACT (last two sub-sections).
$start = $d.Paragraphs(70).Range.Start
$end = $d.Paragraphs(79).Range.End
$d.Range($start, $end).Delete()

# DATA PRODUCTS MANAGER, ANALYTICS SUPERVISOR, SOFTWARE ENGINEER,
# SENIOR ANALYST & PLATFORM DEVELOPER, RESEARCH DIRECTOR & PLATFORM
# ARCHITECT job entries (everything between the first job's bullets
# and the KEY ACHIEVEMENTS heading).
$start = $d.Paragraphs(23).Range.Start
$end = $d.Paragraphs(63).Range.End
$d.Range($start, $end).Delete()

# ------------------------------------------------------------------
# 2. Simple text substitutions.
# ------------------------------------------------------------------

$d.Content.Find.Execute("DHEERAJ CHAND", $true, $true, $false, $false, $false, $true, 1, $false, "Dheeraj Chand", 2)

$d.Content.Find.Execute("Senior Software Engineer & Geospatial Platform Architect", $true, $true, $false, $false, $false, $true, 1, $false, "Professional Title", 2)

$d.Content.Find.Execute("(202) 550-7110 | Dheeraj.Chand@gmail.com", $true, $true, $false, $false, $false, $true, 1, $false, "202.550.7110 | dheeraj.chand@gmail.com", 2)

$d.Content.Find.Execute("Senior Software Engineer with 20+ years building scalable geospatial data platforms, web applications, and distributed analytical systems. Expert in full-stack development with deep specialization in Apache Spark/Sedona for big data geospatial processing. Proven track record architecting multi-tenant SaaS platforms like BALLISTA and DAMON used by thousands of analysts, implementing ETL pipelines processing billions of geospatial records, and building production systems integrating ESRI, OSGeo, and SAFE FME technologies.", $true, $true, $false, $false, $false, $true, 1, $false, "Senior Software Engineer with 21 years building scalable geospatial data platforms, web applications, and distributed analytical systems. Expert in full-stack development with deep specialization in Apache Spark/Sedona for big data geospatial processing. Proven track record architecting multi-tenant SaaS platforms used by thousands of analysts, implementing ETL pipelines processing billions of geospatial records, and building production systems integrating ESRI, OSGeo, and SAFE FME technologies.", 2)

$d.Content.Find.Execute("Siege Analytics, Austin, TX | 2005 – Present", $true, $true, $false, $false, $false, $true, 1, $false, "Your Company Name, Your City, ST | 2005 – Present", 2)

$d.Content.Find.Execute("▸ Architected and engineered BALLISTA: GeoDjango redistricting platform serving thousands of analysts with real-time collaborative editing, Census integration, and legal compliance analysis", $true, $true, $false, $false, $false, $true, 1, $false, "▸ Architected and engineered redistricting platform serving thousands of analysts with real-time collaborative editing, Census integration, and legal compliance analysis", 2)

$d.Content.Find.Execute("▸ Developed DAMON: Flask/PostGIS microservice using incomplete data for boundary estimation without machine learning, processing geographies at national scale", $true, $true, $false, $false, $false, $true, 1, $false, "▸ Developed boundary estimation microservice using incomplete data for boundary estimation without machine learning, processing geographies at national scale", 2)

$d.Content.Find.Execute("✓ Architected BALLISTA redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration", $true, $true, $false, $false, $false, $true, 1, $false, "✓ Architected redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration", 2)

$d.Content.Find.Execute("✓ Built DAMON boundary estimation system achieving accurate geospatial results without machine learning using advanced PostGIS algorithms", $true, $true, $false, $false, $false, $true, 1, $false, "✓ Built boundary estimation system achieving accurate geospatial results without machine learning using advanced PostGIS algorithms", 2)

$d.Content.Find.Execute("✓ Developed SimCrisis econometric simulation platform with NetLogo multi-agent modeling and GeoDjango web interface", $true, $true, $false, $false, $false, $true, 1, $false, "✓ Developed econometric simulation platform with NetLogo multi-agent modeling and web interface", 2)

$d.Content.Find.Execute("✓ Created RACSO comprehensive survey platform managing complete research lifecycle with integrated geospatial market segmentation", $true, $true, $false, $false, $false, $true, 1, $false, "✓ Created comprehensive survey platform managing complete research lifecycle with integrated geospatial market segmentation", 2)
